# Auto-generated edit script: refresh Chocobo Leve profit-tracker market data
# (H/I/J/K/L/M/N columns per leve row), per scheduled-runner update.
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1797.4286
$ws.Range("I6").Value = 194
$ws.Range("K6").Value = 582
$ws.Range("M6").Value = -470
# Row 38
$ws.Range("H38").Value = 2576.7
$ws.Range("I38").Value = 193.8
$ws.Range("J38").Value = 4959.6
$ws.Range("K38").Value = 581.4000000000001
$ws.Range("L38").Value = 14878.8
$ws.Range("M38").Value = -209.4000000000001
$ws.Range("N38").Value = -15622.8
# Row 58
$ws.Range("H58").Value = 21959.8
$ws.Range("J58").Value = 21959.8
$ws.Range("L58").Value = 65879.39999999999
$ws.Range("N58").Value = -66179.39999999999
# Row 87
$ws.Range("H87").Value = 21946.824
$ws.Range("J87").Value = 21946.824
$ws.Range("L87").Value = 21946.824
$ws.Range("N87").Value = -24442.824
# Row 90
$ws.Range("H90").Value = 21946.824
$ws.Range("J90").Value = 21946.824
$ws.Range("L90").Value = 65840.47200000001
$ws.Range("N90").Value = -78320.47200000001
# Row 138
$ws.Range("H138").Value = 3577.0574
$ws.Range("I138").Value = 1477.2273
$ws.Range("J138").Value = 4287.769
$ws.Range("K138").Value = 4431.6819
$ws.Range("L138").Value = 12863.307
$ws.Range("M138").Value = 708.3181000000004
$ws.Range("N138").Value = -23143.307

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 1730
$ws.Range("I102").Value = 1622
$ws.Range("K102").Value = 1622
$ws.Range("M102").Value = 0
# Row 122
$ws.Range("H122").Value = 4081.7273
$ws.Range("I122").Value = 1979.8
$ws.Range("K122").Value = 5939.4
$ws.Range("M122").Value = -3489.4
# Row 132
$ws.Range("H132").Value = 2737.513
$ws.Range("I132").Value = 1470.4584
$ws.Range("K132").Value = 4411.3752
$ws.Range("M132").Value = -1881.3752

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 18520160
$ws.Range("I94").Value = 31251616
$ws.Range("K94").Value = 31251616
$ws.Range("M94").Value = -31251165
# Row 99
$ws.Range("H99").Value = 3848.4736
$ws.Range("I99").Value = 1839.8
$ws.Range("J99").Value = 4565.857
$ws.Range("K99").Value = 1839.8
$ws.Range("L99").Value = 4565.857
$ws.Range("M99").Value = -341.8
$ws.Range("N99").Value = -7561.857
# Row 134
$ws.Range("H134").Value = 3513.0286
$ws.Range("I134").Value = 2170.4138
$ws.Range("K134").Value = 6511.241399999999
$ws.Range("M134").Value = -3976.241399999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7944.4443
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 7944.4443
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 7944.4443
$ws.Range("M31").ClearContents() | Out-Null
$ws.Range("N31").Value = -8534.444299999999
# Row 34
$ws.Range("H34").Value = 7944.4443
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 7944.4443
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 7944.4443
$ws.Range("M34").ClearContents() | Out-Null
$ws.Range("N34").Value = -8348.444299999999
# Row 134
$ws.Range("H134").Value = 20502
$ws.Range("I134").Value = 26003
$ws.Range("J134").Value = 9500
$ws.Range("K134").Value = 78009
$ws.Range("L134").Value = 28500
$ws.Range("M134").Value = -75474
$ws.Range("N134").Value = -33570
# Row 139
$ws.Range("H139").Value = 112540
$ws.Range("J139").Value = 112540
$ws.Range("L139").Value = 112540
$ws.Range("N139").Value = -122820

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 3571514.5
$ws.Range("I2").Value = 89.75
$ws.Range("J2").Value = 5952464.5
$ws.Range("K2").Value = 538.5
$ws.Range("L2").Value = 35714787
$ws.Range("M2").Value = -425.5
$ws.Range("N2").Value = -35715013
# Row 12
$ws.Range("H12").Value = 54.6
$ws.Range("I12").Value = 107.5
$ws.Range("J12").Value = 31.928572
$ws.Range("K12").Value = 322.5
$ws.Range("L12").Value = 95.78571599999999
$ws.Range("M12").Value = -149.5
$ws.Range("N12").Value = -441.785716
# Row 17
$ws.Range("H17").Value = 804.2
$ws.Range("I17").Value = 780
$ws.Range("J17").Value = 901
$ws.Range("K17").Value = 2340
$ws.Range("L17").Value = 2703
$ws.Range("M17").Value = -2171
$ws.Range("N17").Value = -3041
# Row 92
$ws.Range("H92").Value = 41669892
$ws.Range("I92").Value = 736.4
$ws.Range("J92").Value = 71433570
$ws.Range("K92").Value = 2209.2
$ws.Range("L92").Value = 214300710
$ws.Range("M92").Value = -961.1999999999998
$ws.Range("N92").Value = -214303206
# Row 108
$ws.Range("H108").Value = 2222
$ws.Range("I108").Value = 2222
$ws.Range("K108").Value = 6666
$ws.Range("M108").Value = -3786
# Row 110
$ws.Range("H110").Value = 7014
$ws.Range("I110").Value = 4027
$ws.Range("J110").Value = 10001
$ws.Range("K110").Value = 12081
$ws.Range("L110").Value = 30003
$ws.Range("M110").Value = -7991
$ws.Range("N110").Value = -38183
# Row 131
$ws.Range("H131").Value = 9264059
$ws.Range("I131").Value = 45474100
$ws.Range("J131").Value = 1025.0233
$ws.Range("K131").Value = 136422300
$ws.Range("L131").Value = 3075.0699
$ws.Range("M131").Value = -136417260
$ws.Range("N131").Value = -13155.0699

$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
# Row 101
$ws.Range("H101").Value = 48000
$ws.Range("J101").Value = 48000
$ws.Range("L101").Value = 48000
$ws.Range("N101").Value = -54490
# Row 119
$ws.Range("H119").Value = 39800
$ws.Range("J119").Value = 39800
$ws.Range("L119").Value = 39800
$ws.Range("N119").Value = -49476
# Row 125
$ws.Range("H125").Value = 34970
$ws.Range("J125").Value = 34970
$ws.Range("L125").Value = 34970
$ws.Range("N125").Value = -39890
# Row 126
$ws.Range("H126").Value = 4240.94
$ws.Range("I126").Value = 2990
$ws.Range("J126").Value = 5404.6045
$ws.Range("K126").Value = 8970
$ws.Range("L126").Value = 16213.8135
$ws.Range("M126").Value = -6500
$ws.Range("N126").Value = -21153.8135
# Row 132
$ws.Range("H132").Value = 2867.4333
$ws.Range("I132").Value = 731.6923
$ws.Range("J132").Value = 4500.647
$ws.Range("K132").Value = 2195.0769
$ws.Range("L132").Value = 13501.941
$ws.Range("M132").Value = 334.9231
$ws.Range("N132").Value = -18561.941

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2259.1177
$ws.Range("J46").Value = 2158.6667
$ws.Range("L46").Value = 2158.6667
$ws.Range("N46").Value = -2534.6667
# Row 55
$ws.Range("H55").Value = 421.26666
$ws.Range("I55").Value = 349.85715
$ws.Range("K55").Value = 349.85715
$ws.Range("M55").Value = -176.85715
# Row 68
$ws.Range("H68").Value = 943.13336
$ws.Range("I68").Value = 724.64105
$ws.Range("J68").Value = 2363.3333
$ws.Range("K68").Value = 724.64105
$ws.Range("L68").Value = 2363.3333
$ws.Range("M68").Value = 24.35895000000005
$ws.Range("N68").Value = -3861.3333
# Row 71
$ws.Range("H71").Value = 943.13336
$ws.Range("I71").Value = 724.64105
$ws.Range("J71").Value = 2363.3333
$ws.Range("K71").Value = 3623.20525
$ws.Range("L71").Value = 11816.6665
$ws.Range("M71").Value = 120.79475
$ws.Range("N71").Value = -19304.6665
# Row 122
$ws.Range("H122").Value = 3576.7058
$ws.Range("I122").Value = 2500.1304
$ws.Range("J122").Value = 5827.727
$ws.Range("K122").Value = 7500.3912
$ws.Range("L122").Value = 17483.181
$ws.Range("M122").Value = -5050.3912
$ws.Range("N122").Value = -22383.181

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 23865136
$ws.Range("I62").Value = 45457920
$ws.Range("J62").Value = 113072.4
$ws.Range("K62").Value = 45457920
$ws.Range("L62").Value = 113072.4
$ws.Range("M62").Value = -45457296
$ws.Range("N62").Value = -114320.4
# Row 65
$ws.Range("H65").Value = 23865136
$ws.Range("I65").Value = 45457920
$ws.Range("J65").Value = 113072.4
$ws.Range("K65").Value = 227289600
$ws.Range("L65").Value = 565362
$ws.Range("M65").Value = -227286480
$ws.Range("N65").Value = -571602
# Row 81
$ws.Range("H81").Value = 35715652
$ws.Range("I81").Value = 40179860
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 80359720
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -80358659
$ws.Range("N81").Value = -6122
# Row 84
$ws.Range("H84").Value = 35715652
$ws.Range("I84").Value = 40179860
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 401798600
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -401793296
$ws.Range("N84").Value = -30608
# Row 103
$ws.Range("H103").Value = 36200
$ws.Range("J103").Value = 36200
$ws.Range("L103").Value = 36200
$ws.Range("N103").Value = -38544
# Row 122
$ws.Range("H122").Value = 3878.96
$ws.Range("I122").Value = 2233.7646
$ws.Range("J122").Value = 7375
$ws.Range("K122").Value = 6701.293799999999
$ws.Range("L122").Value = 22125
$ws.Range("M122").Value = -4251.293799999999
$ws.Range("N122").Value = -27025
# Row 136
$ws.Range("H136").Value = 3265.8667
$ws.Range("I136").Value = 850.4706
$ws.Range("J136").Value = 6424.4614
$ws.Range("K136").Value = 2551.4118
$ws.Range("L136").Value = 19273.3842
$ws.Range("M136").Value = -1.411799999999857
$ws.Range("N136").Value = -24373.3842
# Row 138
$ws.Range("H138").Value = 54850
$ws.Range("J138").Value = 54850
$ws.Range("L138").Value = 54850
$ws.Range("N138").Value = -65130
